$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Ecological): delete empty column G, shifting H:M left to G:L ---
$ws1.Columns.Item(7).Delete()

# --- Column D width update ---
$ws1.Columns.Item(4).ColumnWidth = 41

# --- New header G1 ---
$ws1.Range("G1").Value = 'Data Source'

# --- Remove bold font from species-name cells in column A (rows below) ---
$ws1.Range("A55").Font.Bold = $false
$ws1.Range("A57").Font.Bold = $false
$ws1.Range("A63").Font.Bold = $false
$ws1.Range("A64").Font.Bold = $false
$ws1.Range("A65").Font.Bold = $false
$ws1.Range("A66").Font.Bold = $false
$ws1.Range("A67").Font.Bold = $false
$ws1.Range("A68").Font.Bold = $false
$ws1.Range("A70").Font.Bold = $false
$ws1.Range("A71").Font.Bold = $false
$ws1.Range("A74").Font.Bold = $false
$ws1.Range("A77").Font.Bold = $false
$ws1.Range("A78").Font.Bold = $false
$ws1.Range("A95").Font.Bold = $false
$ws1.Range("A108").Font.Bold = $false

# --- Populate new "Data Source" column G for Biological trait rows 55-108 ---
$ws1.Range("G55").Value = 'Amniotes; AnAge; NA; Pacifici et al. (2013); PHYLACINE; Smith et al. (2003) (EltonTraits); split from [species name]'
$ws1.Range("G56").Value = 'Tsuboi et al. (2018); imputed; Heldstab et al. (2018); NA; split from [species name]'
$ws1.Range("G57").Value = 'Amniotes; imputed; NA; mean of female and female head body length (Amniotes); female head body length maturity (Amniotes); undefined sex head body length (Amniotes); split from [species name]'
$ws1.Range("G58").Value = 'NA; PanTHERIA; split from [species name]'
$ws1.Range("G59").Value = 'Amniotes; imputed; NA; AnAge; Pacifici et al. (2013); split from [species name]'
$ws1.Range("G60").Value = 'PanTHERIA; NA; mean of female and male (Amniotes); mean of female and undefined sex (Amniotes); female maturity (Amniotes); mean of female, male and undefined (Amniotes); split from [species name]'
$ws1.Range("G61").Value = 'Amniotes; imputed; NA; AnAge; split from [species name]'
$ws1.Range("G62").Value = 'Amniotes; NA; AnAge; split from [species name]'
$ws1.Range("G63").Value = 'Pacifici et al. (2013); calculated; imputed; PanTHERIA; NA; split from [species name]'
$ws1.Range("G64").Value = 'Amniotes; imputed; NA; AnAge; split from [species name]'
$ws1.Range("G65").Value = 'NA; PanTHERIA; split from [species name]'
$ws1.Range("G66").Value = 'Amniotes; imputed; NA; AnAge; split from [species name]'
$ws1.Range("G67").Value = 'Amniotes; imputed; NA; ; AnAge; split from [species name]'
$ws1.Range("G68").Value = 'Amniotes; imputed; NA; AnAge; split from [species name]'
$ws1.Range("G69").Value = 'Amniotes; NA; AnAge; PanTHERIA; split from [species name]'
$ws1.Range("G70").Value = 'Amniotes; imputed; AnAge; NA; split from [species name]'
$ws1.Range("G71").Value = 'Amniotes; NA; PanTHERIA; split from [species name]'
$ws1.Range("G72").Value = 'Pacifici et al. (2013); imputed; IUCN; NA; split from [species name]'
$ws1.Range("G73").Value = 'calculated; imputed; NA; split from [species name]'
$ws1.Range("G74").Value = 'PanTHERIA; NA; TetraDENSITY; split from [species name]'
$ws1.Range("G75").Value = 'Heldstab et al. (2018); imputed; Buckley et al. (2018); Botero et al. (2013); Turbill et al. (2011); marine; NA; split from [species name]'
$ws1.Range("G76").Value = 'imputed; PanTHERIA; NA; split from [species name]'
$ws1.Range("G77").Value = 'NA; PanTHERIA; split from [species name]'
$ws1.Range("G78").Value = 'PanTHERIA; NA; split from [species name]'
$ws1.Range("G79").Value = 'PHYLACINE; NA; imputed; split from [species name]'
$ws1.Range("G80").Value = 'PHYLACINE; NA; imputed; split from [species name]'
$ws1.Range("G81").Value = 'PHYLACINE; NA; imputed; split from [species name]'
$ws1.Range("G82").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G83").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G84").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G85").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G86").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G87").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G88").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G89").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G90").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G91").Value = 'EltonTraits; NA; split from [species name]'
$ws1.Range("G92").Value = 'EltonTraits; imputed; NA; split from [species name]'
$ws1.Range("G93").Value = 'MammalDIET2; PanTHERIA; imputed; NA; split from [species name]'
$ws1.Range("G94").Value = 'EltonTraits; imputed; NA; split from [species name]'
$ws1.Range("G95").Value = 'EltonTraits; PanTHERIA; imputed; NA; split from [species name]'
$ws1.Range("G96").Value = 'IUCN; PHYLACINE; NA; split from [species name]'
$ws1.Range("G97").Value = 'IUCN; PHYLACINE; NA; split from [species name]'
$ws1.Range("G98").Value = 'IUCN; PHYLACINE; NA; split from [species name]'
$ws1.Range("G99").Value = 'IUCN; PHYLACINE; NA; split from [species name]'
$ws1.Range("G100").Value = 'NA; IUCN; split from [species name]'
$ws1.Range("G101").Value = 'NA; IUCN; split from [species name]'
$ws1.Range("G102").Value = 'NA; calculated; split from [species name]'
$ws1.Range("G103").Value = 'Botero et al. (2013); NA; split from [species name]'
$ws1.Range("G104").Value = 'PHYLACINE; NA; split from [species name]'
$ws1.Range("G105").Value = 'Botero et al. (2013); NA; split from [species name]'
$ws1.Range("G106").Value = 'Botero et al. (2013); NA; split from [species name]'
$ws1.Range("G107").Value = 'IUCN; NA; split from [species name]'
$ws1.Range("G108").Value = 'calculated; NA; split from [species name]'

# --- Sheet2 (Biological): add two new trait rows ---
$ws2.Range("A2").Value = 'TNFR_presence'
$ws2.Range("A3").Value = 'TNFR_truncated'

# --- Restore view selections/active sheet ---
$ws2.Range("A8").Select()
$ws1.Select()
$ws1.Range("A70").Select()